# Add team record (Wins / Losses / Ties) columns to the player sheet.
# New columns are appended right after the existing data (through column AC),
# occupying AD:AF, and the sheet dimension grows from A1:AC49 to A1:AF49.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 49

# --- Header row (row 1): new labels, matching the bold/bordered header style
# already used by the other header cells (e.g. AC1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)   # xlPasteFormats

# --- Data rows (2..49): every player row gets the same 1998 team record.
$ws.Range("AD2:AD" + $lastRow).Value = 92
$ws.Range("AE2:AE" + $lastRow).Value = 70
$ws.Range("AF2:AF" + $lastRow).Value = 0
